# Apply changes described in the commit "Cambios en la Demanda"
$wb = $excel.ActiveWorkbook

# --- BY_Data sheet: set zoom to 190% ---
$wsByData = $wb.Worksheets.Item("BY_Data")
$wsByData.Activate()
$excel.ActiveWindow.Zoom = 190

# --- ELC_Dem sheet: set zoom to 145%, update growth rate J4, change selection ---
$wsElcDem = $wb.Worksheets.Item("ELC_Dem")
$wsElcDem.Activate()
$excel.ActiveWindow.Zoom = 145

# Update the annual growth rate from 4% to 3%
$wsElcDem.Range("J4").Value = 0.03

# Update the active selection to F4
$wsElcDem.Range("F4").Select()
